$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text: Volume number 37 -> 38 ---
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "38"

# --- Update header text: report week dates ---
$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 8).Text = "9/15/2025"
$weekCell.Characters(47, 9).Text = "9/21/2025"

# --- Update crime-statistics data cells (rows 14-31) ---
$ws.Range("M14").Value = -80
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -33.333333333333
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = -16.666666666666
$ws.Range("I15").Value = 58
$ws.Range("J15").Value = 47
$ws.Range("K15").Value = 23.404255319148
$ws.Range("L15").Value = 107.142857142857
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = -10.769230769230
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = 4.761904761904
$ws.Range("I16").Value = 219
$ws.Range("J16").Value = 222
$ws.Range("K16").Value = -1.351351351351
$ws.Range("L16").Value = -7.594936708860
$ws.Range("M16").Value = -26.755852842809
$ws.Range("N16").Value = -77.258566978193
$ws.Range("C17").Value = 24
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 84
$ws.Range("H17").Value = 21.739130434782
$ws.Range("I17").Value = 754
$ws.Range("J17").Value = 702
$ws.Range("K17").Value = 7.407407407407
$ws.Range("L17").Value = 6.047819971870
$ws.Range("M17").Value = 112.994350282486
$ws.Range("N17").Value = -14.898419864559
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -37.5
$ws.Range("F18").Value = 23
$ws.Range("H18").Value = -34.285714285714
$ws.Range("I18").Value = 222
$ws.Range("J18").Value = 242
$ws.Range("K18").Value = -8.264462809917
$ws.Range("L18").Value = -12.252964426877
$ws.Range("M18").Value = -50
$ws.Range("N18").Value = -91.141260973663
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = -34.375
$ws.Range("F19").Value = 123
$ws.Range("G19").Value = 112
$ws.Range("H19").Value = 9.821428571428
$ws.Range("I19").Value = 1062
$ws.Range("J19").Value = 1128
$ws.Range("K19").Value = -5.851063829787
$ws.Range("L19").Value = -8.841201716738
$ws.Range("M19").Value = 51.931330472103
$ws.Range("N19").Value = -9.308283518360
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -27.586206896551
$ws.Range("I20").Value = 166
$ws.Range("J20").Value = 257
$ws.Range("K20").Value = -35.408560311284
$ws.Range("L20").Value = -49.848942598187
$ws.Range("M20").Value = -35.408560311284
$ws.Range("N20").Value = -95.265259555048
$ws.Range("C21").Value = 62
$ws.Range("D21").Value = 78
$ws.Range("E21").Value = -20.512820512820
$ws.Range("F21").Value = 278
$ws.Range("G21").Value = 272
$ws.Range("H21").Value = 2.205882352941
$ws.Range("I21").Value = 2483
$ws.Range("J21").Value = 2603
$ws.Range("K21").Value = -4.610065309258
$ws.Range("L21").Value = -9.412623130244
$ws.Range("M21").Value = 18.069424631478
$ws.Range("N21").Value = -72.759188151398
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = -53.333333333333
$ws.Range("I23").Value = 108
$ws.Range("J23").Value = 96
$ws.Range("K23").Value = 12.5
$ws.Range("L23").Value = 4.854368932038
$ws.Range("M23").Value = 107.692307692308
$ws.Range("C24").Value = 53
$ws.Range("D24").Value = 67
$ws.Range("E24").Value = -20.895522388059
$ws.Range("F24").Value = 254
$ws.Range("G24").Value = 279
$ws.Range("H24").Value = -8.960573476702
$ws.Range("I24").Value = 2933
$ws.Range("J24").Value = 2924
$ws.Range("K24").Value = 0.307797537619
$ws.Range("L24").Value = -4.865390853065
$ws.Range("M24").Value = 6.191165821868
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 34
$ws.Range("E25").Value = -29.411764705882
$ws.Range("F25").Value = 127
$ws.Range("G25").Value = 149
$ws.Range("H25").Value = -14.765100671140
$ws.Range("I25").Value = 1703
$ws.Range("J25").Value = 1551
$ws.Range("K25").Value = 9.800128949065
$ws.Range("L25").Value = 18.263888888888
$ws.Range("C26").Value = 47
$ws.Range("D26").Value = 40
$ws.Range("E26").Value = 17.5
$ws.Range("F26").Value = 170
$ws.Range("G26").Value = 159
$ws.Range("H26").Value = 6.918238993710
$ws.Range("I26").Value = 1404
$ws.Range("J26").Value = 1387
$ws.Range("K26").Value = 1.225666906993
$ws.Range("L26").Value = 8.500772797527
$ws.Range("M26").Value = -5.454545454545
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -22.222222222222
$ws.Range("I27").Value = 73
$ws.Range("J27").Value = 77
$ws.Range("K27").Value = -5.194805194805
$ws.Range("L27").Value = 52.083333333333
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 20
$ws.Range("H28").Value = -30
$ws.Range("I28").Value = 156
$ws.Range("J28").Value = 149
$ws.Range("K28").Value = 4.697986577181
$ws.Range("L28").Value = 0
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("C29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("D29").Value = 2
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E29").Value = -50
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -33.333333333333
$ws.Range("I29").Value = 8
$ws.Range("J29").Value = 13
$ws.Range("K29").Value = -38.461538461538
$ws.Range("L29").Value = -66.666666666666
$ws.Range("M29").Value = -66.666666666666
$ws.Range("N29").Value = -90.123456790123
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("C30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 2
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -33.333333333333
$ws.Range("I30").Value = 8
$ws.Range("J30").Value = 13
$ws.Range("K30").Value = -38.461538461538
$ws.Range("L30").Value = -63.636363636363
$ws.Range("M30").Value = -63.636363636363
$ws.Range("N30").Value = -88.405797101449
$ws.Range("D31").Value = 2
$ws.Range("J31").Value = 19
$ws.Range("K31").Value = -42.105263157894
$ws.Range("L31").Value = -8.333333333333
